$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Content.Find.Execute("2024-10-13 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-14 Monday", 2) | Out-Null

# Update each arithmetic-expression cell in the table, addressed by
# (row, column) so the one duplicated source string ("39-16=") still
# maps to the correct distinct replacement per position.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "88-69="
$t.Cell(1, 2).Range.Text = "8+25="
$t.Cell(1, 3).Range.Text = "92-56="
$t.Cell(1, 4).Range.Text = "42-0="
$t.Cell(1, 5).Range.Text = "68-11="

$t.Cell(2, 1).Range.Text = "65+25="
$t.Cell(2, 2).Range.Text = "67-9="
$t.Cell(2, 3).Range.Text = "57-41="
$t.Cell(2, 4).Range.Text = "47+15="
$t.Cell(2, 5).Range.Text = "63+19="

$t.Cell(3, 1).Range.Text = "7+63="
$t.Cell(3, 2).Range.Text = "77-44="
$t.Cell(3, 3).Range.Text = "33+7="
$t.Cell(3, 4).Range.Text = "41+2="
$t.Cell(3, 5).Range.Text = "73+22="

$t.Cell(4, 1).Range.Text = "77-28="
$t.Cell(4, 2).Range.Text = "95-32="
$t.Cell(4, 3).Range.Text = "20-7="
$t.Cell(4, 4).Range.Text = "56-28="
$t.Cell(4, 5).Range.Text = "28+61="

$t.Cell(5, 1).Range.Text = "48-31="
$t.Cell(5, 2).Range.Text = "77-13="
$t.Cell(5, 3).Range.Text = "35+42="
$t.Cell(5, 4).Range.Text = "99-81="
$t.Cell(5, 5).Range.Text = "39-25="

$t.Cell(6, 1).Range.Text = "72+17="
$t.Cell(6, 2).Range.Text = "70-68="
$t.Cell(6, 3).Range.Text = "8+35="
$t.Cell(6, 4).Range.Text = "84-14="
$t.Cell(6, 5).Range.Text = "21-8="

$t.Cell(7, 1).Range.Text = "70+27="
$t.Cell(7, 2).Range.Text = "72-64="
$t.Cell(7, 3).Range.Text = "24+3="
$t.Cell(7, 4).Range.Text = "90-63="
$t.Cell(7, 5).Range.Text = "13+35="

$t.Cell(8, 1).Range.Text = "72-69="
$t.Cell(8, 2).Range.Text = "91-86="
$t.Cell(8, 3).Range.Text = "35+37="
$t.Cell(8, 4).Range.Text = "28+38="
$t.Cell(8, 5).Range.Text = "60+11="

$t.Cell(9, 1).Range.Text = "59+5="
$t.Cell(9, 2).Range.Text = "2+2="
$t.Cell(9, 3).Range.Text = "0+79="
$t.Cell(9, 4).Range.Text = "31+55="
$t.Cell(9, 5).Range.Text = "12+10="

$t.Cell(10, 1).Range.Text = "6+38="
$t.Cell(10, 2).Range.Text = "57-12="
$t.Cell(10, 3).Range.Text = "76-55="
$t.Cell(10, 4).Range.Text = "96-92="
$t.Cell(10, 5).Range.Text = "14+36="

$t.Cell(11, 1).Range.Text = "91-40="
$t.Cell(11, 2).Range.Text = "74-10="
$t.Cell(11, 3).Range.Text = "98-44="
$t.Cell(11, 4).Range.Text = "2+74="
$t.Cell(11, 5).Range.Text = "64+25="

$t.Cell(12, 1).Range.Text = "82-0="
$t.Cell(12, 2).Range.Text = "57-39="
$t.Cell(12, 3).Range.Text = "0+87="
$t.Cell(12, 4).Range.Text = "59-38="
$t.Cell(12, 5).Range.Text = "21-1="

$t.Cell(13, 1).Range.Text = "59+3="
$t.Cell(13, 2).Range.Text = "32+28="
$t.Cell(13, 3).Range.Text = "18+39="
$t.Cell(13, 4).Range.Text = "20+52="
$t.Cell(13, 5).Range.Text = "61+7="

$t.Cell(14, 1).Range.Text = "16+78="
$t.Cell(14, 2).Range.Text = "81-52="
$t.Cell(14, 3).Range.Text = "92-1="
$t.Cell(14, 4).Range.Text = "60-55="
$t.Cell(14, 5).Range.Text = "59+20="

$t.Cell(15, 1).Range.Text = "54-40="
$t.Cell(15, 2).Range.Text = "33-32="
$t.Cell(15, 3).Range.Text = "0+90="
$t.Cell(15, 4).Range.Text = "92-24="
$t.Cell(15, 5).Range.Text = "94-51="

$t.Cell(16, 1).Range.Text = "96-25="
$t.Cell(16, 2).Range.Text = "25+65="
$t.Cell(16, 3).Range.Text = "44+8="
$t.Cell(16, 4).Range.Text = "78+0="
$t.Cell(16, 5).Range.Text = "59-38="

$t.Cell(17, 1).Range.Text = "70-19="
$t.Cell(17, 2).Range.Text = "74-22="
$t.Cell(17, 3).Range.Text = "93-48="
$t.Cell(17, 4).Range.Text = "38+46="
$t.Cell(17, 5).Range.Text = "14-12="

$t.Cell(18, 1).Range.Text = "37-14="
$t.Cell(18, 2).Range.Text = "85-21="
$t.Cell(18, 3).Range.Text = "46+46="
$t.Cell(18, 4).Range.Text = "81-24="
$t.Cell(18, 5).Range.Text = "23-2="

$t.Cell(19, 1).Range.Text = "13+80="
$t.Cell(19, 2).Range.Text = "86-60="
$t.Cell(19, 3).Range.Text = "98-74="
$t.Cell(19, 4).Range.Text = "90+2="
$t.Cell(19, 5).Range.Text = "26-13="

$t.Cell(20, 1).Range.Text = "42-4="
$t.Cell(20, 2).Range.Text = "66-52="
$t.Cell(20, 3).Range.Text = "86-26="
$t.Cell(20, 4).Range.Text = "12+0="
$t.Cell(20, 5).Range.Text = "9+32="

